# Realestate Update resale numbers 2023-06-18 14:07
# Append a new data row (row 55) to the CityResaleNum sheet with the
# latest resale numbers snapshot, matching the existing columns:
# A Date, B Time, C Weekday, D Week (text), E..T city values (numeric).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 55

# Columns A-D hold text values in this sheet (dates/times/weekday/week are
# stored as text, not as real Excel dates/numbers), so force text format
# before assigning so Excel doesn't auto-convert "2023-06-18" into a date
# serial number or "25" into a number.
$textRange = $ws.Range("A$row`:D$row")
$textRange.NumberFormat = "@"

$ws.Range("A$row").Value = "2023-06-18"
$ws.Range("B$row").Value = "14:04:21"
$ws.Range("C$row").Value = "Sunday"
$ws.Range("D$row").Value = "25"

# Columns E-T hold numeric resale counts per city.
$ws.Range("E$row").Value = 122080
$ws.Range("F$row").Value = 133587
$ws.Range("G$row").Value = 162271
$ws.Range("H$row").Value = 133155
$ws.Range("I$row").Value = 177357
$ws.Range("J$row").Value = 114917
$ws.Range("K$row").Value = 201346
$ws.Range("L$row").Value = 225016
$ws.Range("M$row").Value = 175377
$ws.Range("N$row").Value = 103639
$ws.Range("O$row").Value = 39146
$ws.Range("P$row").Value = 33980
$ws.Range("Q$row").Value = 51789
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 36411
$ws.Range("T$row").Value = -1
